# Apply "repull data, push all data, mean calculation" edits:
# Update column F (dSF) values for a subset of rows to reflect newly
# recalculated/repulled data that now differs from column E (dS0).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -3
    3  = -5
    5  = 3
    10 = 0
    17 = 0
    19 = -5
    20 = 2
    44 = 2
    47 = -1
    50 = 0
    51 = 1
    53 = 0
    64 = -5
    65 = 0
    66 = -2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
